$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the exposure-site details on row 2 with the new Black Rock entry
$ws.Range("A2").Value = "Black Rock"
$ws.Range("B2").Value = "Smile Buffalo Thai restaurant  305 Beach Road, Black Rock VIC 3193"
$ws.Range("C2").Value = "21/12/20 7:30pm-9:30pm"
$ws.Range("D2").Value = "Case dined in restaurant"
$ws.Range("E2").Value = "new"

# The old (now superseded) entry on row 3 is removed entirely
$ws.Rows.Item(3).Delete()

# Match the resulting selection state
$ws.Range("B2").Select()
